$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 49.5
$ws.Cells.Item(5, 9).Value = 49.5
$ws.Cells.Item(5, 11).Value = 49.5
$ws.Cells.Item(5, 13).Value = 65.5

$ws.Cells.Item(86, 8).Value = 66698904
$ws.Cells.Item(86, 9).Value = 2525.889
$ws.Cells.Item(86, 10).Value = 166743470
$ws.Cells.Item(86, 11).Value = 2525.889
$ws.Cells.Item(86, 12).Value = 166743470
$ws.Cells.Item(86, 13).Value = -1402.889
$ws.Cells.Item(86, 14).Value = -166745716

$ws.Cells.Item(89, 8).Value = 66698904
$ws.Cells.Item(89, 9).Value = 2525.889
$ws.Cells.Item(89, 10).Value = 166743470
$ws.Cells.Item(89, 11).Value = 12629.445
$ws.Cells.Item(89, 12).Value = 833717350
$ws.Cells.Item(89, 13).Value = -7013.445
$ws.Cells.Item(89, 14).Value = -833728582

$ws.Cells.Item(112, 8).Value = 3665
$ws.Cells.Item(112, 9).Value = 1191.75
$ws.Cells.Item(112, 10).Value = 4077.2083
$ws.Cells.Item(112, 11).Value = 3575.25
$ws.Cells.Item(112, 12).Value = 12231.6249
$ws.Cells.Item(112, 13).Value = -2467.25
$ws.Cells.Item(112, 14).Value = -14447.6249

$ws.Cells.Item(132, 8).Value = 5640.9414
$ws.Cells.Item(132, 9).Value = 5865.5713
$ws.Cells.Item(132, 11).Value = 17596.7139
$ws.Cells.Item(132, 13).Value = -15066.7139

$ws.Cells.Item(135, 8).Value = 3521.6
$ws.Cells.Item(135, 9).Value = 3861.0454
$ws.Cells.Item(135, 11).Value = 34749.4086
$ws.Cells.Item(135, 13).Value = -32214.4086

$ws.Cells.Item(138, 8).Value = 427831.7
$ws.Cells.Item(138, 9).Value = 879582.75
$ws.Cells.Item(138, 10).Value = 4315.0625
$ws.Cells.Item(138, 11).Value = 2638748.25
$ws.Cells.Item(138, 12).Value = 12945.1875
$ws.Cells.Item(138, 13).Value = -2633608.25
$ws.Cells.Item(138, 14).Value = -23225.1875

$ws.Cells.Item(141, 8).Value = 5156.5
$ws.Cells.Item(141, 9).Value = 4783.923
$ws.Cells.Item(141, 11).Value = 14351.769
$ws.Cells.Item(141, 13).Value = -9171.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3510.34
$ws.Cells.Item(61, 9).Value = 3396.0454
$ws.Cells.Item(61, 11).Value = 3396.0454
$ws.Cells.Item(61, 13).Value = -3184.0454

$ws.Cells.Item(97, 8).Value = 8701148
$ws.Cells.Item(97, 9).Value = 6971
$ws.Cells.Item(97, 10).Value = 33334650
$ws.Cells.Item(97, 11).Value = 6971
$ws.Cells.Item(97, 12).Value = 33334650
$ws.Cells.Item(97, 13).Value = -6475
$ws.Cells.Item(97, 14).Value = -33335642

$ws.Cells.Item(110, 8).Value = 2816.625
$ws.Cells.Item(110, 9).Value = 1008.25
$ws.Cells.Item(110, 11).Value = 1008.25
$ws.Cells.Item(110, 13).Value = 1036.75

$ws.Cells.Item(136, 8).Value = 3510.34
$ws.Cells.Item(136, 9).Value = 3396.0454
$ws.Cells.Item(136, 11).Value = 10188.1362
$ws.Cells.Item(136, 13).Value = -7638.136200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 21784.273
$ws.Cells.Item(82, 9).Value = 3537
$ws.Cells.Item(82, 10).Value = 53717
$ws.Cells.Item(82, 11).Value = 3537
$ws.Cells.Item(82, 12).Value = 53717
$ws.Cells.Item(82, 13).Value = -3154
$ws.Cells.Item(82, 14).Value = -54483

$ws.Cells.Item(85, 8).Value = 21784.273
$ws.Cells.Item(85, 9).Value = 3537
$ws.Cells.Item(85, 10).Value = 53717
$ws.Cells.Item(85, 11).Value = 3537
$ws.Cells.Item(85, 12).Value = 53717
$ws.Cells.Item(85, 13).Value = -2211
$ws.Cells.Item(85, 14).Value = -56369

$ws.Cells.Item(86, 8).Value = 5416.522
$ws.Cells.Item(86, 9).Value = 6745.625
$ws.Cells.Item(86, 10).Value = 2378.5715
$ws.Cells.Item(86, 11).Value = 6745.625
$ws.Cells.Item(86, 12).Value = 2378.5715
$ws.Cells.Item(86, 13).Value = -5622.625
$ws.Cells.Item(86, 14).Value = -4624.5715

$ws.Cells.Item(89, 8).Value = 5416.522
$ws.Cells.Item(89, 9).Value = 6745.625
$ws.Cells.Item(89, 10).Value = 2378.5715
$ws.Cells.Item(89, 11).Value = 33728.125
$ws.Cells.Item(89, 12).Value = 11892.8575
$ws.Cells.Item(89, 13).Value = -28112.125
$ws.Cells.Item(89, 14).Value = -23124.8575

$ws.Cells.Item(94, 8).Value = 10308.387
$ws.Cells.Item(94, 9).Value = 12092.375
$ws.Cells.Item(94, 11).Value = 12092.375
$ws.Cells.Item(94, 13).Value = -11641.375

$ws.Cells.Item(107, 8).Value = 1601.4706
$ws.Cells.Item(107, 9).Value = 1482.0667
$ws.Cells.Item(107, 10).Value = 2497
$ws.Cells.Item(107, 11).Value = 1482.0667
$ws.Cells.Item(107, 12).Value = 2497
$ws.Cells.Item(107, 13).Value = 437.9332999999999
$ws.Cells.Item(107, 14).Value = -6337

$ws.Cells.Item(134, 8).Value = 7037.909
$ws.Cells.Item(134, 9).Value = 7563.41
$ws.Cells.Item(134, 11).Value = 22690.23
$ws.Cells.Item(134, 13).Value = -20155.23

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2192.5
$ws.Cells.Item(58, 9).Value = 1869
$ws.Cells.Item(58, 10).Value = 3163
$ws.Cells.Item(58, 11).Value = 1869
$ws.Cells.Item(58, 12).Value = 3163
$ws.Cells.Item(58, 13).Value = -1666
$ws.Cells.Item(58, 14).Value = -3569

$ws.Cells.Item(62, 8).Value = 50196.785
$ws.Cells.Item(62, 9).Value = 8225
$ws.Cells.Item(62, 10).Value = 81675.625
$ws.Cells.Item(62, 11).Value = 8225
$ws.Cells.Item(62, 12).Value = 81675.625
$ws.Cells.Item(62, 13).Value = -7601
$ws.Cells.Item(62, 14).Value = -82923.625

$ws.Cells.Item(65, 8).Value = 50196.785
$ws.Cells.Item(65, 9).Value = 8225
$ws.Cells.Item(65, 10).Value = 81675.625
$ws.Cells.Item(65, 11).Value = 41125
$ws.Cells.Item(65, 12).Value = 408378.125
$ws.Cells.Item(65, 13).Value = -38005
$ws.Cells.Item(65, 14).Value = -414618.125

$ws.Cells.Item(122, 8).Value = 10161.571
$ws.Cells.Item(122, 9).Value = 25235.2
$ws.Cells.Item(122, 11).Value = 75705.60000000001
$ws.Cells.Item(122, 13).Value = -73255.60000000001

$ws.Cells.Item(136, 8).Value = 2192.5
$ws.Cells.Item(136, 9).Value = 1869
$ws.Cells.Item(136, 10).Value = 3163
$ws.Cells.Item(136, 11).Value = 5607
$ws.Cells.Item(136, 12).Value = 9489
$ws.Cells.Item(136, 13).Value = -3057
$ws.Cells.Item(136, 14).Value = -14589

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 827.6667
$ws.Cells.Item(13, 9).Value = 820.5
$ws.Cells.Item(13, 11).Value = 2461.5
$ws.Cells.Item(13, 13).Value = -2293.5

$ws.Cells.Item(107, 8).Value = 361.07693
$ws.Cells.Item(107, 9).Value = 311.75
$ws.Cells.Item(107, 10).Value = 440
$ws.Cells.Item(107, 11).Value = 935.25
$ws.Cells.Item(107, 12).Value = 1320
$ws.Cells.Item(107, 13).Value = 984.75
$ws.Cells.Item(107, 14).Value = -5160

$ws.Cells.Item(121, 8).Value = 1720830.6
$ws.Cells.Item(121, 9).Value = 1808867.9
$ws.Cells.Item(121, 10).Value = 1669475.6
$ws.Cells.Item(121, 11).Value = 5426603.699999999
$ws.Cells.Item(121, 12).Value = 5008426.800000001
$ws.Cells.Item(121, 13).Value = -5425293.699999999
$ws.Cells.Item(121, 14).Value = -5011046.800000001

$ws.Cells.Item(129, 8).Value = 17546056
$ws.Cells.Item(129, 9).Value = 736
$ws.Cells.Item(129, 10).Value = 30306288
$ws.Cells.Item(129, 11).Value = 2208
$ws.Cells.Item(129, 12).Value = 90918864
$ws.Cells.Item(129, 13).Value = 2792
$ws.Cells.Item(129, 14).Value = -90928864

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6192.75
$ws.Cells.Item(70, 9).Value = 5427.5
$ws.Cells.Item(70, 10).Value = 7723.25
$ws.Cells.Item(70, 11).Value = 5427.5
$ws.Cells.Item(70, 12).Value = 7723.25
$ws.Cells.Item(70, 13).Value = -5157.5
$ws.Cells.Item(70, 14).Value = -8263.25

$ws.Cells.Item(73, 8).Value = 6192.75
$ws.Cells.Item(73, 9).Value = 5427.5
$ws.Cells.Item(73, 10).Value = 7723.25
$ws.Cells.Item(73, 11).Value = 5427.5
$ws.Cells.Item(73, 12).Value = 7723.25
$ws.Cells.Item(73, 13).Value = -4491.5
$ws.Cells.Item(73, 14).Value = -9595.25

$ws.Cells.Item(102, 8).Value = 6822.346
$ws.Cells.Item(102, 9).Value = 9726.4375
$ws.Cells.Item(102, 11).Value = 9726.4375
$ws.Cells.Item(102, 13).Value = -8104.4375

$ws.Cells.Item(107, 8).Value = 395.5
$ws.Cells.Item(107, 10).Value = 80
$ws.Cells.Item(107, 12).Value = 80
$ws.Cells.Item(107, 14).Value = -3920

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 5081.467
$ws.Cells.Item(16, 9).Value = 5202.75
$ws.Cells.Item(16, 11).Value = 5202.75
$ws.Cells.Item(16, 13).Value = -5032.75

$ws.Cells.Item(40, 8).Value = 19775.758
$ws.Cells.Item(40, 9).Value = 31174.5
$ws.Cells.Item(40, 10).Value = 9136.933999999999
$ws.Cells.Item(40, 11).Value = 31174.5
$ws.Cells.Item(40, 12).Value = 9136.933999999999
$ws.Cells.Item(40, 13).Value = -31038.5
$ws.Cells.Item(40, 14).Value = -9408.933999999999

$ws.Cells.Item(82, 8).Value = 2794.56
$ws.Cells.Item(82, 9).Value = 2669.5
$ws.Cells.Item(82, 10).Value = 2953.7273
$ws.Cells.Item(82, 11).Value = 2669.5
$ws.Cells.Item(82, 12).Value = 2953.7273
$ws.Cells.Item(82, 13).Value = -2308.5
$ws.Cells.Item(82, 14).Value = -3675.7273

$ws.Cells.Item(85, 8).Value = 2794.56
$ws.Cells.Item(85, 9).Value = 2669.5
$ws.Cells.Item(85, 10).Value = 2953.7273
$ws.Cells.Item(85, 11).Value = 2669.5
$ws.Cells.Item(85, 12).Value = 2953.7273
$ws.Cells.Item(85, 13).Value = -1421.5
$ws.Cells.Item(85, 14).Value = -5449.7273

$ws.Cells.Item(132, 8).Value = 881722.9
$ws.Cells.Item(132, 9).Value = 1867037.8
$ws.Cells.Item(132, 10).Value = 5887.4443
$ws.Cells.Item(132, 11).Value = 5601113.4
$ws.Cells.Item(132, 12).Value = 17662.3329
$ws.Cells.Item(132, 13).Value = -5598583.4
$ws.Cells.Item(132, 14).Value = -22722.3329

$ws.Cells.Item(136, 8).Value = 4590.486
$ws.Cells.Item(136, 9).Value = 2480.1667
$ws.Cells.Item(136, 10).Value = 9194.817999999999
$ws.Cells.Item(136, 11).Value = 7440.500100000001
$ws.Cells.Item(136, 12).Value = 27584.454
$ws.Cells.Item(136, 13).Value = -4890.500100000001
$ws.Cells.Item(136, 14).Value = -32684.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 43777.785
$ws.Cells.Item(100, 10).Value = 113333
$ws.Cells.Item(100, 12).Value = 226666
$ws.Cells.Item(100, 14).Value = -227748

$ws.Cells.Item(101, 8).Value = 35000
$ws.Cells.Item(101, 10).Value = 35000
$ws.Cells.Item(101, 12).Value = 35000
$ws.Cells.Item(101, 14).Value = -41490

$ws.Cells.Item(103, 8).Value = 38400
$ws.Cells.Item(103, 10).Value = 38400
$ws.Cells.Item(103, 12).Value = 38400
$ws.Cells.Item(103, 14).Value = -40744

$ws.Cells.Item(105, 8).Value = 45750
$ws.Cells.Item(105, 10).Value = 45750
$ws.Cells.Item(105, 12).Value = 45750
$ws.Cells.Item(105, 14).Value = -52738

$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 1650.2693
$ws.Cells.Item(113, 9).Value = 800.55554
$ws.Cells.Item(113, 11).Value = 2401.66662
$ws.Cells.Item(113, 13).Value = -231.66662

$ws.Cells.Item(126, 8).Value = 16008.733
$ws.Cells.Item(126, 9).Value = 20779.38
$ws.Cells.Item(126, 10).Value = 4877.222
$ws.Cells.Item(126, 11).Value = 62338.14
$ws.Cells.Item(126, 12).Value = 14631.666
$ws.Cells.Item(126, 13).Value = -59868.14
$ws.Cells.Item(126, 14).Value = -19571.666

$ws.Cells.Item(132, 8).Value = 9874.574000000001
$ws.Cells.Item(132, 9).Value = 11974.383
$ws.Cells.Item(132, 10).Value = 4382.769
$ws.Cells.Item(132, 11).Value = 35923.149
$ws.Cells.Item(132, 12).Value = 13148.307
$ws.Cells.Item(132, 13).Value = -33393.149
$ws.Cells.Item(132, 14).Value = -18208.307

$ws.Cells.Item(136, 8).Value = 322560.22
$ws.Cells.Item(136, 9).Value = 482764.6
$ws.Cells.Item(136, 10).Value = 2151.5
$ws.Cells.Item(136, 11).Value = 1448293.8
$ws.Cells.Item(136, 12).Value = 6454.5
$ws.Cells.Item(136, 13).Value = -1445743.8
$ws.Cells.Item(136, 14).Value = -11554.5
